# Swap switch for if block.
# Adds an "accept / reject / clip" breakdown table (rows 27-31) below the
# existing benchmark table, plus a small scratch rounding-error check in
# K24:K26, and nudges the current selection to H19.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- small floating point sanity-check block (K24:K26) ---------------------
$ws.Range("K24").Value = 8.747
$ws.Range("K25").Value = 8.746
$ws.Range("K26").Formula = "=K24-K25"

# --- new accept / reject / clip breakdown table (rows 27-31) ---------------
$ws.Range("B27").Value = "accept"
$ws.Range("C27").Value = 305986
$ws.Range("D27").Formula = "=C27/C31"
$ws.Range("D27").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("F27").Value = 3995
$ws.Range("F27").NumberFormat = "0"
$ws.Range("G27").Formula = "=F27/F31"
$ws.Range("G27").NumberFormat = $ws.Range("D22").NumberFormat

$ws.Range("B28").Value = "reject"
$ws.Range("C28").Value = 498925
$ws.Range("D28").Formula = "=C28/C31"
$ws.Range("D28").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("F28").Value = 0
$ws.Range("F28").NumberFormat = "0"
$ws.Range("G28").Formula = "=F28/F31"
$ws.Range("G28").NumberFormat = $ws.Range("D22").NumberFormat

$ws.Range("B29").Value = "clip"
$ws.Range("C29").Formula = "=813434-(C27+C28)"
$ws.Range("D29").Formula = "=C29/C31"
$ws.Range("D29").NumberFormat = $ws.Range("D22").NumberFormat
$ws.Range("F29").Formula = "=6427-F27"
$ws.Range("F29").NumberFormat = "0"
$ws.Range("G29").Formula = "=F29/F31"
$ws.Range("G29").NumberFormat = $ws.Range("D22").NumberFormat

$ws.Range("F30").NumberFormat = "0"

$ws.Range("C31").Formula = "=SUM(C27:C30)"
$ws.Range("F31").Formula = "=SUM(F27:F30)"
$ws.Range("F31").NumberFormat = "0"

# --- move the saved selection, matching the author's last cursor spot ------
$ws.Range("H19").Select()

$wb.Save()
